$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'89.315.58"
$ws.Range("E2").Formula = "'  -1.58%  "
$ws.Range("D3").Formula = "'3.088.17"
$ws.Range("E3").Formula = "'  -2.50%  "
$ws.Range("D4").Formula = "'1.00"
$ws.Range("E4").Formula = "'  +0.10%  "
$ws.Range("D5").Formula = "'212.71"
$ws.Range("E5").Formula = "'  -1.47%  "
$ws.Range("D6").Formula = "'622.01"
$ws.Range("E6").Formula = "'  -0.84%  "
$ws.Range("D7").Formula = "'0.370"
$ws.Range("E7").Formula = "'  -6.35%  "
$ws.Range("D8").Formula = "'0.821"
$ws.Range("E8").Formula = "'  +16.63%  "
$ws.Range("D9").Formula = "'1.00"
$ws.Range("E9").Formula = "'  +0.14%  "
$ws.Range("D10").Formula = "'3.084.74"
$ws.Range("E10").Formula = "'  -2.56%  "
$ws.Range("D11").Formula = "'0.612"
$ws.Range("E11").Formula = "'  +9.14%  "
$ws.Range("E12").Formula = "'  +0.76%  "
$ws.Range("D13").Formula = "'0.0000240"
$ws.Range("E13").Formula = "'  -4.74%  "
$ws.Range("D14").Formula = "'5.29"
$ws.Range("E14").Formula = "'  -0.34%  "
$ws.Range("D15").Formula = "'89.084.65"
$ws.Range("E15").Formula = "'  -1.56%  "
$ws.Range("D16").Formula = "'32.18"
$ws.Range("E16").Formula = "'  -0.55%  "
$ws.Range("D17").Formula = "'3.665.34"
$ws.Range("E17").Formula = "'  -2.16%  "
$ws.Range("D18").Formula = "'3.097.51"
$ws.Range("E18").Formula = "'  -3.11%  "
$ws.Range("D19").Formula = "'3.38"
$ws.Range("E19").Formula = "'  +2.86%  "
$ws.Range("E20").Formula = "'  +1.13%  "
$ws.Range("D21").Formula = "'13.47"
$ws.Range("E21").Formula = "'  +2.01%  "
$ws.Range("D22").Formula = "'424.37"
$ws.Range("E22").Formula = "'  -2.50%  "
$ws.Range("D23").Formula = "'8.27"
$ws.Range("E23").Formula = "'  -1.79%  "
$ws.Range("D24").Formula = "'4.94"
$ws.Range("E24").Formula = "'  -0.29%  "
$ws.Range("D25").Formula = "'5.52"
$ws.Range("E25").Formula = "'  +7.33%  "
$ws.Range("D26").Formula = "'12.09"
$ws.Range("E26").Formula = "'  +4.28%  "
$ws.Range("D27").Formula = "'83.27"
$ws.Range("E27").Formula = "'  +4.26%  "
$ws.Range("D28").Formula = "'3.257.42"
$ws.Range("E28").Formula = "'  -2.54%  "
$ws.Range("D29").Formula = "'1.00"
$ws.Range("E29").Formula = "'  +0.16%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Formula = "'1.07"
$ws.Range("E30").Formula = "'  +7.16%  "
$ws.Range("B31").Value = "Cronos"
$ws.Range("C31").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D31").Formula = "'0.169"
$ws.Range("E31").Formula = "'  +9.15%  "
$ws.Range("D32").Formula = "'8.09"
$ws.Range("E32").Formula = "'  -1.74%  "
$ws.Range("D33").Formula = "'506.33"
$ws.Range("E33").Formula = "'  -2.49%  "
$ws.Range("E34").Formula = "'  -8.98%  "
$ws.Range("D35").Formula = "'6.70"
$ws.Range("E35").Formula = "'  -2.61%  "
$ws.Range("D36").Formula = "'1.26"
$ws.Range("E36").Formula = "'  -1.06%  "
$ws.Range("E37").Formula = "'  -4.00%  "
$ws.Range("E38").Formula = "'  +0.69%  "
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").Formula = "'22.29"
$ws.Range("E39").Formula = "'  -0.40%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Formula = "'0.128"
$ws.Range("E40").Formula = "'  +3.10%  "
$ws.Range("E41").Formula = "'  +0.25%  "
$ws.Range("E42").Formula = "'  +0.01%  "
$ws.Range("D43").Formula = "'0.363"
$ws.Range("E43").Formula = "'  -0.69%  "
$ws.Range("E44").Formula = "'  -3.86%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Formula = "'0.134"
$ws.Range("E45").Formula = "'  +7.92%  "
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Formula = "'144.96"
$ws.Range("E46").Formula = "'  -0.27%  "
$ws.Range("E47").Formula = "'  +16.43%  "
$ws.Range("D48").Formula = "'43.40"
$ws.Range("E48").Formula = "'  -1.56%  "
$ws.Range("B49").Value = "ImmutableX"
$ws.Range("C49").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D49").Formula = "'1.23"
$ws.Range("E49").Formula = "'  +1.82%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Formula = "'159.72"
$ws.Range("E50").Formula = "'  -6.07%  "
$ws.Range("D51").Formula = "'0.705"
$ws.Range("E51").Formula = "'  -4.83%  "
